$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05045533333333333
$ws.Range("H2").Value = 0.151366
$ws.Range("I2").Value = 0.004442474524580737
$ws.Range("J2").Value = 0.004442474524580737
$ws.Range("M2").Value = 2.440259666666666
$ws.Range("N2").Value = 7.320779
$ws.Range("O2").Value = 0.6045788173334784
$ws.Range("P2").Value = 0.6045788173334784
$ws.Range("Q2").Value = 0.1231241149015555
$ws.Range("R2").Value = 1.108117034114
$ws.Range("S2").Value = 0.002685825994105129
$ws.Range("T2").Value = 0.002685825994105129

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05045533333333333
$ws.Range("H3").Value = 0.151366
$ws.Range("I3").Value = 0.004442474524580737
$ws.Range("J3").Value = 0.004442474524580737
$ws.Range("O3").Value = 0.1795692107559644
$ws.Range("P3").Value = 0.1795692107559644
$ws.Range("Q3").Value = 0.03656975650488888
$ws.Range("R3").Value = 0.329127808544
$ws.Range("S3").Value = 0.000797731644182441
$ws.Range("T3").Value = 0.000797731644182441

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05045533333333333
$ws.Range("H4").Value = 0.151366
$ws.Range("I4").Value = 0.004442474524580737
$ws.Range("J4").Value = 0.004442474524580737
$ws.Range("M4").Value = 0.8712426666666667
$ws.Range("N4").Value = 2.613728
$ws.Range("O4").Value = 0.2158519719105573
$ws.Range("P4").Value = 0.2158519719105573
$ws.Range("Q4").Value = 0.04395883916088889
$ws.Range("R4").Value = 0.395629552448
$ws.Range("S4").Value = 0.0009589168862931678
$ws.Range("T4").Value = 0.0009589168862931677

# Row 5
$ws.Range("I5").Value = 0.7425623198471305
$ws.Range("J5").Value = 0.7425623198471305
$ws.Range("M5").Value = 2.440259666666666
$ws.Range("N5").Value = 7.320779
$ws.Range("O5").Value = 0.6045788173334784
$ws.Range("P5").Value = 0.6045788173334784
$ws.Range("Q5").Value = 20.58027072176678
$ws.Range("R5").Value = 185.222436495901
$ws.Range("S5").Value = 0.4489374491295823
$ws.Range("T5").Value = 0.4489374491295823

# Row 6
$ws.Range("I6").Value = 0.7425623198471305
$ws.Range("J6").Value = 0.7425623198471305
$ws.Range("O6").Value = 0.1795692107559644
$ws.Range("P6").Value = 0.1795692107559644
$ws.Range("S6").Value = 0.1333413297120672
$ws.Range("T6").Value = 0.1333413297120672

# Row 7
$ws.Range("I7").Value = 0.7425623198471305
$ws.Range("J7").Value = 0.7425623198471305
$ws.Range("M7").Value = 0.8712426666666667
$ws.Range("N7").Value = 2.613728
$ws.Range("O7").Value = 0.2158519719105573
$ws.Range("P7").Value = 0.2158519719105573
$ws.Range("Q7").Value = 7.347746712892445
$ws.Range("R7").Value = 66.129720416032
$ws.Range("S7").Value = 0.1602835410054811
$ws.Range("T7").Value = 0.1602835410054811

# Row 8
$ws.Range("G8").Value = 2.873389
$ws.Range("H8").Value = 8.620167
$ws.Range("I8").Value = 0.2529952056282888
$ws.Range("J8").Value = 0.2529952056282888
$ws.Range("M8").Value = 2.440259666666666
$ws.Range("N8").Value = 7.320779
$ws.Range("O8").Value = 0.6045788173334784
$ws.Range("P8").Value = 0.6045788173334784
$ws.Range("Q8").Value = 7.011815283343666
$ws.Range("R8").Value = 63.106337550093
$ws.Range("S8").Value = 0.152955542209791
$ws.Range("T8").Value = 0.152955542209791

# Row 9
$ws.Range("G9").Value = 2.873389
$ws.Range("H9").Value = 8.620167
$ws.Range("I9").Value = 0.2529952056282888
$ws.Range("J9").Value = 0.2529952056282888
$ws.Range("O9").Value = 0.1795692107559644
$ws.Range("P9").Value = 0.1795692107559644
$ws.Range("Q9").Value = 2.082617022458666
$ws.Range("R9").Value = 18.743553202128
$ws.Range("S9").Value = 0.04543014939971473
$ws.Range("T9").Value = 0.04543014939971473

# Row 10
$ws.Range("G10").Value = 2.873389
$ws.Range("H10").Value = 8.620167
$ws.Range("I10").Value = 0.2529952056282888
$ws.Range("J10").Value = 0.2529952056282888
$ws.Range("M10").Value = 0.8712426666666667
$ws.Range("N10").Value = 2.613728
$ws.Range("O10").Value = 0.2158519719105573
$ws.Range("P10").Value = 0.2158519719105573
$ws.Range("Q10").Value = 2.503419094730667
$ws.Range("R10").Value = 22.530771852576
$ws.Range("S10").Value = 0.05460951401878307
$ws.Range("T10").Value = 0.05460951401878306
